$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.941.53"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "3.586.29"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("D5").Value = "578.94"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").Value = "190.72"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("D8").Value = "3.582.75"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  +1.96%  "
$ws.Range("D11").Value = "0.666"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "55.74"
$ws.Range("E12").Value = "  -4.01%  "
$ws.Range("D13").Value = "0.0000306"
$ws.Range("E13").Value = "  +5.35%  "
$ws.Range("D14").Value = "9.65"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").Value = "4.162.49"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").Value = "3.582.55"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "69.919.61"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").Value = "474.70"
$ws.Range("E22").Value = "  -3.78%  "
$ws.Range("D23").Value = "19.37"
$ws.Range("E23").Value = "  +12.61%  "
$ws.Range("D24").Value = "5.03"
$ws.Range("E24").Value = "  -6.10%  "
$ws.Range("D25").Value = "95.56"
$ws.Range("E25").Value = "  +5.27%  "
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("E27").Value = "  -3.49%  "
$ws.Range("D28").Value = "11.03"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").Value = "9.29"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").Value = "32.25"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").Value = "12.21"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").Value = "66.44"
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("D35").Value = "590.34"
$ws.Range("D36").Value = "39.01"
$ws.Range("E36").Value = "  +2.47%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").Value = "0.0₃0801"
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("E39").Value = "  -2.18%  "
$ws.Range("D40").Value = "3.19"
$ws.Range("E40").Value = "  +17.74%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.137"
$ws.Range("E41").Value = "  -5.94%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "3.46"
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("D43").Value = "3.223.04"
$ws.Range("E43").Value = "  -2.35%  "
$ws.Range("E44").Value = "  +6.73%  "
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "3.35"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("D48").Value = "9.44"
$ws.Range("E48").Value = "  +3.45%  "
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  -5.32%  "
